$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.262.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.46%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.792.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.63%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'315.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.17%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.35%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.5399"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.49%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3761"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.55%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.30%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'41.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.09%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'1.094"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.41%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.24%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'20.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.40%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.090"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.77%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "'Chainlink"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'7.271"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.95%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'WrappedEther"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1.786.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.68%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'89.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.76%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.00001055"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.24%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06518"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.02%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'17.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.01%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.937"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.70%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'28.283.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.42%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.40%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.084"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.25%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'158.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.69%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'20.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.68%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.989.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.05%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.285"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.70%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'121.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.78%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.086"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.23%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.1046"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.52%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'3.663"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.09%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.535"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.39%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.2259"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.83%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.06465"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.46%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.02278"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.79%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -0.53%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'8.487"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.66%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.6147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.43%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.187"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.62%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'WEMIXTOKEN"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.440"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.28%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'11.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.26%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.22%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'13.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.38%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'3.671"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.28%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5756"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.57%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'125.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.42%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.188"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.57%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.927"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.14%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.06847"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.84%  "
$ws.Range("E51").Style = "Normal"
